# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for Mango (Vega Central Mapocho de
# Santiago) at row 548, pushing the existing rows 548-572 down to 549-573.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(548).Insert()

$ws.Cells.Item(548, 1).Value = 9
$ws.Cells.Item(548, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(548, 3).Value = "Metropolitana"
$ws.Cells.Item(548, 4).Value = 44939
$ws.Cells.Item(548, 5).Value = 13
$ws.Cells.Item(548, 6).Value = "Fruta"
$ws.Cells.Item(548, 7).Value = 100108
$ws.Cells.Item(548, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(548, 9).Value = 100108002
$ws.Cells.Item(548, 10).Value = "Mango"
$ws.Cells.Item(548, 11).Value = "Sin especificar"
$ws.Cells.Item(548, 12).Value = "Primera"
$ws.Cells.Item(548, 13).Value = 630
$ws.Cells.Item(548, 14).Value = 5500
$ws.Cells.Item(548, 15).Value = 6000
$ws.Cells.Item(548, 16).Value = 5722
$ws.Cells.Item(548, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(548, 18).Value = "Perú"
$ws.Cells.Item(548, 19).Value = 1430
$ws.Cells.Item(548, 20).Value = 4
